# Se implementa la opcion de actualizar producto
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "clientes" to "productos"
$ws.Name = "productos"

# Renumber existing product IDs (keep them as text, matching column formatting)
$ws.Range("A2").Value = "1201"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1202"

# Add the new product row (plain/unstyled cells, like the rest of row 3)
$ws.Range("A4").Value = "1203"
$ws.Range("B4").Value = "iPhone 16 Pro"
$ws.Range("C4").Value = "999"
$ws.Range("D4").Value = "5"
$ws.Range("A4:D4").Style = "Normal"

$null = $ws.Range("A4").Select()
